# Add two new header/value columns (I, J) to the sheet, matching the
# existing header style used by the neighboring "IP" column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell H1 onto the new
# header cells I1:J1 so they pick up the same style (bold font,
# border, centered alignment) without creating a duplicate style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
